$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text representation
# (values like "213.98" or "3.00" would otherwise be auto-coerced to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.114.83"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.623.45"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "213.98"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").Value = "0.522"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "20.29"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "1.632.97"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "0.542"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "27.095.61"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "64.53"
$ws.Range("E16").Value = "  -4.27%  "
$ws.Range("D17").Value = "0.0₃0744"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "216.32"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "6.92"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("E22").Value = "  -6.86%  "
$ws.Range("D23").Value = "9.04"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").Value = "148.07"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "7.28"
$ws.Range("E26").Value = "  -3.55%  "
$ws.Range("D27").Value = "0.117"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "15.59"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").Value = "3.35"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "3.00"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "1.343.67"
$ws.Range("E33").Value = "  +5.54%  "
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Value = "0.858"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "0.803"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "65.60"
$ws.Range("E41").Value = "  +6.00%  "
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").Value = "1.760.80"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "0.906"
$ws.Range("E45").Value = "  +35.39%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "90.58"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "0.0995"
$ws.Range("E50").Value = "  +1.89%  "
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").Value = "  -0.58%  "
